# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 235 (pushing the existing rows
# 235..396 down to 236..397) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(235).Insert()

$ws.Range("A235").Value = 3
$ws.Range("B235").Value = "Femacal de La Calera"
$ws.Range("C235").Value = "Coquimbo"
$ws.Range("D235").Value = 44957
$ws.Range("D235").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E235").Value = 5
$ws.Range("F235").Value = 100112039
$ws.Range("G235").Value = "Ciboulette"
$ws.Range("H235").Value = "Sin especificar"
$ws.Range("I235").Value = "Primera"
$ws.Range("J235").Value = 120
$ws.Range("K235").Value = 1500
$ws.Range("L235").Value = 1500
$ws.Range("M235").Value = 1500
$ws.Range("N235").Value = "`$/docena de atados"
$ws.Range("O235").Value = "Provincia de Quillota"
$ws.Range("P235").Value = 500
$ws.Range("Q235").Value = 3
$ws.Range("R235").Value = "Hortaliza"
